$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update as a quote-prefixed literal so Excel stores it as
# text (matching the workbooks existing inlineStr/text cells) instead of
# re-interpreting dotted price strings (e.g. "1.003", "27.492.82") as numbers.
function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.Value = "'" + $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "27.492.82"
Set-TextValue "E2" "  +5.08%  "
Set-TextValue "D3" "1.724.07"
Set-TextValue "E3" "  +4.05%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "226.01"
Set-TextValue "E5" "  +3.23%  "
Set-TextValue "D6" "0.5383"
Set-TextValue "E6" "  +2.67%  "
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "D8" "0.2683"
Set-TextValue "E8" "  +0.69%  "
Set-TextValue "D9" "0.06605"
Set-TextValue "E9" "  +3.88%  "
Set-TextValue "D10" "21.75"
Set-TextValue "E10" "  +5.62%  "
Set-TextValue "D11" "0.07742"
Set-TextValue "E11" "  +0.69%  "
Set-TextValue "D12" "4.640"
Set-TextValue "E12" "  +0.49%  "
Set-TextValue "D13" "1.730.82"
Set-TextValue "E13" "  +4.40%  "
Set-TextValue "D14" "1.961.30"
Set-TextValue "D15" "0.5885"
Set-TextValue "E15" "  +4.63%  "
Set-TextValue "D16" "0.0₅8286"
Set-TextValue "E16" "  +1.01%  "
Set-TextValue "D17" "68.09"
Set-TextValue "E17" "  +3.83%  "
Set-TextValue "D18" "27.515.91"
Set-TextValue "E18" "  +5.24%  "
Set-TextValue "D19" "222.70"
Set-TextValue "E19" "  +15.03%  "
Set-TextValue "D20" "1.003"
Set-TextValue "E20" "  +0.00%  "
Set-TextValue "D21" "4.741"
Set-TextValue "E21" "  +1.58%  "
Set-TextValue "E22" "  +1.14%  "
Set-TextValue "E23" "  +2.42%  "
Set-TextValue "D25" "147.89"
Set-TextValue "E25" "  +1.74%  "
Set-TextValue "B26" "Toncoin"
Set-TextValue "C26" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D26" "1.693"
Set-TextValue "E26" "  +11.52%  "
Set-TextValue "B27" "Stellar"
Set-TextValue "C27" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D27" "0.1233"
Set-TextValue "E27" "  +2.71%  "
Set-TextValue "D28" "7.423"
Set-TextValue "E28" "  +2.10%  "
Set-TextValue "D29" "16.70"
Set-TextValue "E29" "  +4.46%  "
Set-TextValue "E30" "  +1.30%  "
Set-TextValue "E31" "  +2.54%  "
Set-TextValue "D32" "3.551"
Set-TextValue "E32" "  +2.33%  "
Set-TextValue "D33" "3.467"
Set-TextValue "E33" "  +2.94%  "
Set-TextValue "D34" "1.663"
Set-TextValue "E34" "  +6.34%  "
Set-TextValue "D35" "0.9599"
Set-TextValue "E35" "  +0.67%  "
Set-TextValue "E36" "  +1.84%  "
Set-TextValue "D37" "2.816"
Set-TextValue "E37" "  +1.35%  "
Set-TextValue "D38" "0.5939"
Set-TextValue "E38" "  +4.23%  "
Set-TextValue "D39" "0.01647"
Set-TextValue "E39" "  +3.49%  "
Set-TextValue "D40" "5.867"
Set-TextValue "E40" "  -0.19%  "
Set-TextValue "B41" "Maker"
Set-TextValue "C41" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D41" "1.057.05"
Set-TextValue "E41" "  +2.93%  "
Set-TextValue "B42" "TrustWalletToken"
Set-TextValue "C42" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D42" "0.8546"
Set-TextValue "E42" "  +3.03%  "
Set-TextValue "D43" "1.003"
Set-TextValue "D45" "1.867.35"
Set-TextValue "E46" "  +9.48%  "
Set-TextValue "D47" "58.98"
Set-TextValue "E47" "  +1.61%  "
Set-TextValue "D48" "8.221"
Set-TextValue "E48" "  +2.30%  "
Set-TextValue "D49" "0.4438"
Set-TextValue "E49" "  +2.10%  "
Set-TextValue "D50" "0.9997"
Set-TextValue "E50" "  -0.16%  "
Set-TextValue "D51" "0.05269"
Set-TextValue "E51" "  +1.28%  "
